$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Anpep -> Sele, ECs -> sCs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.710646666666667
$ws.Range("H2").Value = 29.13194
$ws.Range("I2").Value = 0.1041051736296835
$ws.Range("J2").Value = 0.1041051736296835
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.321929333333333
$ws.Range("N2").Value = 21.965788
$ws.Range("Q2").Value = 71.10066867430223
$ws.Range("R2").Value = 639.90601806872
$ws.Range("S2").Value = 0.1041051736296835
$ws.Range("T2").Value = 0.1041051736296835

# Row 3 (FAPs -> sCs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 77.72230533333334
$ws.Range("H3").Value = 233.166916
$ws.Range("I3").Value = 0.8332394710025435
$ws.Range("J3").Value = 0.8332394710025435
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.321929333333333
$ws.Range("N3").Value = 21.965788
$ws.Range("Q3").Value = 569.0772272744232
$ws.Range("R3").Value = 5121.695045469808
$ws.Range("S3").Value = 0.8332394710025435
$ws.Range("T3").Value = 0.8332394710025435

# Row 4 (sCs -> sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.844320666666666
$ws.Range("H4").Value = 17.532962
$ws.Range("I4").Value = 0.06265535536777307
$ws.Range("J4").Value = 0.06265535536777307
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.321929333333333
$ws.Range("N4").Value = 21.965788
$ws.Range("Q4").Value = 42.79170292267288
$ws.Range("R4").Value = 385.1253263040559
$ws.Range("S4").Value = 0.06265535536777307
$ws.Range("T4").Value = 0.06265535536777307
